$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.58102533333333
$ws.Range("H2").Value = 58.743076
$ws.Range("I2").Value = 0.3529199051285138
$ws.Range("J2").Value = 0.3529199051285138
$ws.Range("M2").Value = 5.483544666666667
$ws.Range("N2").Value = 16.450634
$ws.Range("O2").Value = 0.4099269772514247
$ws.Range("P2").Value = 0.4099269772514247
$ws.Range("Q2").Value = 107.3734270344649
$ws.Range("R2").Value = 966.3608433101841
$ws.Range("S2").Value = 0.1446713899211912
$ws.Range("T2").Value = 0.1446713899211912
$ws.Range("G3").Value = 19.58102533333333
$ws.Range("H3").Value = 58.743076
$ws.Range("I3").Value = 0.3529199051285138
$ws.Range("J3").Value = 0.3529199051285138
$ws.Range("O3").Value = 0.32690571008764
$ws.Range("P3").Value = 0.32690571008764
$ws.Range("Q3").Value = 85.62741258113466
$ws.Range("R3").Value = 770.646713230212
$ws.Range("S3").Value = 0.1153715321900993
$ws.Range("T3").Value = 0.1153715321900994
$ws.Range("G4").Value = 19.58102533333333
$ws.Range("H4").Value = 58.743076
$ws.Range("I4").Value = 0.3529199051285138
$ws.Range("J4").Value = 0.3529199051285138
$ws.Range("M4").Value = 2.750327
$ws.Range("N4").Value = 8.250980999999999
$ws.Range("O4").Value = 0.2056029999019452
$ws.Range("P4").Value = 0.2056029999019452
$ws.Range("Q4").Value = 53.85422266195066
$ws.Range("R4").Value = 484.688003957556
$ws.Range("S4").Value = 0.07256139121953235
$ws.Range("T4").Value = 0.07256139121953235
$ws.Range("G5").Value = 19.58102533333333
$ws.Range("H5").Value = 58.743076
$ws.Range("I5").Value = 0.3529199051285138
$ws.Range("J5").Value = 0.3529199051285138
$ws.Range("M5").Value = 0.770031
$ws.Range("N5").Value = 2.310093
$ws.Range("O5").Value = 0.05756431275899004
$ws.Range("P5").Value = 0.05756431275899004
$ws.Range("Q5").Value = 15.077996518452
$ws.Range("R5").Value = 135.701968666068
$ws.Range("S5").Value = 0.02031559179769086
$ws.Range("T5").Value = 0.02031559179769087
$ws.Range("I6").Value = 0.2647616806631773
$ws.Range("J6").Value = 0.2647616806631773
$ws.Range("M6").Value = 5.483544666666667
$ws.Range("N6").Value = 16.450634
$ws.Range("O6").Value = 0.4099269772514247
$ws.Range("P6").Value = 0.4099269772514247
$ws.Range("Q6").Value = 80.55190026716666
$ws.Range("R6").Value = 724.9671024045
$ws.Range("S6").Value = 0.1085329554462632
$ws.Range("T6").Value = 0.1085329554462633
$ws.Range("I7").Value = 0.2647616806631773
$ws.Range("J7").Value = 0.2647616806631773
$ws.Range("O7").Value = 0.32690571008764
$ws.Range("P7").Value = 0.32690571008764
$ws.Range("R7").Value = 578.1417143872499
$ws.Range("S7").Value = 0.08655210522119294
$ws.Range("T7").Value = 0.08655210522119297
$ws.Range("I8").Value = 0.2647616806631773
$ws.Range("J8").Value = 0.2647616806631773
$ws.Range("M8").Value = 2.750327
$ws.Range("N8").Value = 8.250980999999999
$ws.Range("O8").Value = 0.2056029999019452
$ws.Range("P8").Value = 0.2056029999019452
$ws.Range("Q8").Value = 40.40161604824999
$ws.Range("R8").Value = 363.6145444342499
$ws.Range("S8").Value = 0.0544357958034301
$ws.Range("T8").Value = 0.05443579580343011
$ws.Range("I9").Value = 0.2647616806631773
$ws.Range("J9").Value = 0.2647616806631773
$ws.Range("M9").Value = 0.770031
$ws.Range("N9").Value = 2.310093
$ws.Range("O9").Value = 0.05756431275899004
$ws.Range("P9").Value = 0.05756431275899004
$ws.Range("Q9").Value = 11.31156288225
$ws.Range("R9").Value = 101.80406594025
$ws.Range("S9").Value = 0.01524082419229098
$ws.Range("T9").Value = 0.01524082419229099
$ws.Range("G10").Value = 19.14352733333333
$ws.Range("H10").Value = 57.430582
$ws.Range("I10").Value = 0.3450346309906436
$ws.Range("J10").Value = 0.3450346309906436
$ws.Range("M10").Value = 5.483544666666667
$ws.Range("N10").Value = 16.450634
$ws.Range("O10").Value = 0.4099269772514247
$ws.Range("P10").Value = 0.4099269772514247
$ws.Range("Q10").Value = 104.9743872098876
$ws.Range("R10").Value = 944.769484888988
$ws.Range("S10").Value = 0.1414390033290553
$ws.Range("T10").Value = 0.1414390033290553
$ws.Range("G11").Value = 19.14352733333333
$ws.Range("H11").Value = 57.430582
$ws.Range("I11").Value = 0.3450346309906436
$ws.Range("J11").Value = 0.3450346309906436
$ws.Range("O11").Value = 0.32690571008764
$ws.Range("P11").Value = 0.32690571008764
$ws.Range("Q11").Value = 83.71424301459267
$ws.Range("R11").Value = 753.428187131334
$ws.Range("S11").Value = 0.1127937910488232
$ws.Range("T11").Value = 0.1127937910488232
$ws.Range("G12").Value = 19.14352733333333
$ws.Range("H12").Value = 57.430582
$ws.Range("I12").Value = 0.3450346309906436
$ws.Range("J12").Value = 0.3450346309906436
$ws.Range("M12").Value = 2.750327
$ws.Range("N12").Value = 8.250980999999999
$ws.Range("O12").Value = 0.2056029999019452
$ws.Range("P12").Value = 0.2056029999019452
$ws.Range("Q12").Value = 52.65096010010467
$ws.Range("R12").Value = 473.858640900942
$ws.Range("S12").Value = 0.07094015520173702
$ws.Range("T12").Value = 0.07094015520173702
$ws.Range("G13").Value = 19.14352733333333
$ws.Range("H13").Value = 57.430582
$ws.Range("I13").Value = 0.3450346309906436
$ws.Range("J13").Value = 0.3450346309906436
$ws.Range("M13").Value = 0.770031
$ws.Range("N13").Value = 2.310093
$ws.Range("O13").Value = 0.05756431275899004
$ws.Range("P13").Value = 0.05756431275899004
$ws.Range("Q13").Value = 14.741109496014
$ws.Range("R13").Value = 132.669985464126
$ws.Range("S13").Value = 0.01986168141102813
$ws.Range("T13").Value = 0.01986168141102813
$ws.Range("G14").Value = 2.068613
$ws.Range("H14").Value = 6.205839
$ws.Range("I14").Value = 0.03728378321766519
$ws.Range("J14").Value = 0.0372837832176652
$ws.Range("M14").Value = 5.483544666666667
$ws.Range("N14").Value = 16.450634
$ws.Range("O14").Value = 0.4099269772514247
$ws.Range("P14").Value = 0.4099269772514247
$ws.Range("Q14").Value = 11.34333178354733
$ws.Range("R14").Value = 102.089986051926
$ws.Range("S14").Value = 0.01528362855491489
$ws.Range("T14").Value = 0.01528362855491489
$ws.Range("G15").Value = 2.068613
$ws.Range("H15").Value = 6.205839
$ws.Range("I15").Value = 0.03728378321766519
$ws.Range("J15").Value = 0.0372837832176652
$ws.Range("O15").Value = 0.32690571008764
$ws.Range("P15").Value = 0.32690571008764
$ws.Range("Q15").Value = 9.046001208127
$ws.Range("R15").Value = 81.414010873143
$ws.Range("S15").Value = 0.01218828162752447
$ws.Range("T15").Value = 0.01218828162752448
$ws.Range("G16").Value = 2.068613
$ws.Range("H16").Value = 6.205839
$ws.Range("I16").Value = 0.03728378321766519
$ws.Range("J16").Value = 0.0372837832176652
$ws.Range("M16").Value = 2.750327
$ws.Range("N16").Value = 8.250980999999999
$ws.Range("O16").Value = 0.2056029999019452
$ws.Range("P16").Value = 0.2056029999019452
$ws.Range("Q16").Value = 5.689362186451
$ws.Range("R16").Value = 51.204259678059
$ws.Range("S16").Value = 0.007665657677245765
$ws.Range("T16").Value = 0.007665657677245766
$ws.Range("G17").Value = 2.068613
$ws.Range("H17").Value = 6.205839
$ws.Range("I17").Value = 0.03728378321766519
$ws.Range("J17").Value = 0.0372837832176652
$ws.Range("M17").Value = 0.770031
$ws.Range("N17").Value = 2.310093
$ws.Range("O17").Value = 0.05756431275899004
$ws.Range("P17").Value = 0.05756431275899004
$ws.Range("Q17").Value = 1.592896137003
$ws.Range("R17").Value = 14.336065233027
$ws.Range("S17").Value = 0.002146215357980063
$ws.Range("T17").Value = 0.002146215357980063
